$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift data up: drop the 2020 row, update years/values for 2021-2025,
# and remove the now-empty trailing row 7.
$ws.Range("A2").Value = 2021
$ws.Range("B2").Value = 950468.77

$ws.Range("A3").Value = 2022
$ws.Range("B3").Value = 1766113.68

$ws.Range("A4").Value = 2023
$ws.Range("B4").Value = 2842827.64

$ws.Range("A5").Value = 2024
$ws.Range("B5").Value = 4442894.22

$ws.Range("A6").Value = 2025
$ws.Range("B6").Value = 1014612.85

# Remove the old last row entirely so the used range shrinks to A1:B6
$ws.Rows.Item(7).Delete()
